$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.955187333333333
$ws.Range("H2").Value = 8.865562
$ws.Range("I2").Value = 0.0151950616208272
$ws.Range("J2").Value = 0.0151950616208272
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 50.73031201328666
$ws.Range("R2").Value = 456.57280811958
$ws.Range("S2").Value = 0.0008514482838311087
$ws.Range("T2").Value = 0.0008514482838311088

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.955187333333333
$ws.Range("H3").Value = 8.865562
$ws.Range("I3").Value = 0.0151950616208272
$ws.Range("J3").Value = 0.0151950616208272
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 757.8409628265324
$ws.Range("R3").Value = 6820.568665438792
$ws.Range("S3").Value = 0.01271946419423887
$ws.Range("T3").Value = 0.01271946419423888

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.955187333333333
$ws.Range("H4").Value = 8.865562
$ws.Range("I4").Value = 0.0151950616208272
$ws.Range("J4").Value = 0.0151950616208272
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 96.76875781281021
$ws.Range("R4").Value = 870.918820315292
$ws.Range("S4").Value = 0.001624149142757212
$ws.Range("T4").Value = 0.001624149142757213

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 167.1121116666667
$ws.Range("H5").Value = 501.336335
$ws.Range("I5").Value = 0.8592615451885246
$ws.Range("J5").Value = 0.8592615451885246
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 2868.735078289183
$ws.Range("R5").Value = 25818.61570460265
$ws.Range("S5").Value = 0.04814832517757224
$ws.Range("T5").Value = 0.04814832517757224

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 167.1121116666667
$ws.Range("H6").Value = 501.336335
$ws.Range("I6").Value = 0.8592615451885246
$ws.Range("J6").Value = 0.8592615451885246
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 42854.94938914476
$ws.Range("R6").Value = 385694.5445023028
$ws.Range("S6").Value = 0.7192696370860014
$ws.Range("T6").Value = 0.7192696370860014

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 167.1121116666667
$ws.Range("H7").Value = 501.336335
$ws.Range("I7").Value = 0.8592615451885246
$ws.Range("J7").Value = 0.8592615451885246
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 5472.1510474324
$ws.Range("R7").Value = 49249.35942689161
$ws.Range("S7").Value = 0.09184358292495079
$ws.Range("T7").Value = 0.09184358292495079

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 24.41610666666667
$ws.Range("H8").Value = 73.24832
$ws.Range("I8").Value = 0.1255433931906482
$ws.Range("J8").Value = 0.1255433931906482
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 419.1398275765333
$ws.Range("R8").Value = 3772.2584481888
$ws.Range("S8").Value = 0.007034766251424542
$ws.Range("T8").Value = 0.007034766251424543

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 24.41610666666667
$ws.Range("H9").Value = 73.24832
$ws.Range("I9").Value = 0.1255433931906482
$ws.Range("J9").Value = 0.1255433931906482
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 6261.371513077902
$ws.Range("R9").Value = 56352.34361770112
$ws.Range("S9").Value = 0.1050897149586401
$ws.Range("T9").Value = 0.1050897149586401

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 24.41610666666667
$ws.Range("H10").Value = 73.24832
$ws.Range("I10").Value = 0.1255433931906482
$ws.Range("J10").Value = 0.1255433931906482
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 799.514902526791
$ws.Range("R10").Value = 7195.63412274112
$ws.Range("S10").Value = 0.01341891198058352
$ws.Range("T10").Value = 0.01341891198058352

